# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$statusNew = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet updates (status column for each locale) ---
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew

# --- zh-cn sheet updates ---
$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("K2").Value = "2016-08-18 04:50:18"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet updates ---
$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("K2").Value = "2016-08-18 04:50:25"
$wsDeDe.Range("P2").Value = ""

# --- Column width updates ---
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527

$wsZhCn.Range("C1").ColumnWidth = 29.9777047293527
$wsZhCn.Range("P1").ColumnWidth = 13.7470528738839

$wsDeDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDeDe.Range("P1").ColumnWidth = 13.7470528738839
